$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.422.72"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").Value = "1.806.77"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  +0.15%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "227.89"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.56%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.579"
$c.ClearFormats()
$ws.Range("E6").Value = "  +4.07%  "

$ws.Range("E7").Value = "  +0.14%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "35.88"
$c.ClearFormats()
$ws.Range("E8").Value = "  +8.50%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.302"
$c.ClearFormats()
$ws.Range("E9").Value = "  +2.32%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0695"
$c.ClearFormats()
$ws.Range("E10").Value = "  +0.77%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0964"
$c.ClearFormats()
$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("D12").Value = "2.071.11"
$ws.Range("E12").Value = "  +1.24%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "11.47"
$c.ClearFormats()
$ws.Range("E13").Value = "  +2.23%  "

$ws.Range("D14").Value = "1.793.78"
$ws.Range("E14").Value = "  +0.39%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.646"
$c.ClearFormats()
$ws.Range("E15").Value = "  +1.61%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "4.52"
$c.ClearFormats()
$ws.Range("E16").Value = "  +5.16%  "

$ws.Range("D17").Value = "34.416.57"
$ws.Range("E17").Value = "  +0.10%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "69.19"
$c.ClearFormats()
$ws.Range("E18").Value = "  +0.91%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "246.29"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("D20").Value = "0.0₃0798"
$ws.Range("E20").Value = "  -0.14%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.51"
$c.ClearFormats()
$ws.Range("E21").Value = "  +1.67%  "

$ws.Range("E22").Value = "  +0.11%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.20"
$c.ClearFormats()
$ws.Range("E23").Value = "  +0.69%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.14"
$c.ClearFormats()
$ws.Range("E24").Value = "  +4.16%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "171.23"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.88%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.93"
$c.ClearFormats()
$ws.Range("E26").Value = "  +7.83%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "16.98"
$c.ClearFormats()
$ws.Range("E27").Value = "  +2.52%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.119"
$c.ClearFormats()
$ws.Range("E28").Value = "  +3.41%  "

$ws.Range("E29").Value = "  +0.17%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.05"
$c.ClearFormats()
$ws.Range("E30").Value = "  +0.01%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0533"
$c.ClearFormats()
$ws.Range("E31").Value = "  +1.28%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.87"
$c.ClearFormats()
$ws.Range("E32").Value = "  +1.73%  "

$ws.Range("E33").Value = "  +0.49%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.84"
$c.ClearFormats()
$ws.Range("E34").Value = "  +0.97%  "

$ws.Range("D35").Value = "1.399.63"
$ws.Range("E35").Value = "  -0.92%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.674"
$c.ClearFormats()
$ws.Range("E36").Value = "  -1.63%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.48"
$c.ClearFormats()
$ws.Range("E37").Value = "  -3.44%  "

$ws.Range("E38").Value = "  +0.85%  "

$ws.Range("E39").Value = "  +0.46%  "

$ws.Range("E40").Value = "  +13.08%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "82.79"
$c.ClearFormats()
$ws.Range("E41").Value = "  -1.89%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.961"
$c.ClearFormats()
$ws.Range("E42").Value = "  +1.89%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.82"
$c.ClearFormats()
$ws.Range("E43").Value = "  +1.80%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.39"
$c.ClearFormats()
$ws.Range("E45").Value = "  -4.49%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "6.04"
$c.ClearFormats()
$ws.Range("E46").Value = "  -1.07%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0505"
$c.ClearFormats()
$ws.Range("E47").Value = "  -4.03%  "

$ws.Range("D48").Value = "1.969.05"
$ws.Range("E48").Value = "  +1.07%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "105.29"
$c.ClearFormats()
$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("E51").Value = "  +1.53%  "
